# Applies the FineMapping_summary.docx edit: expands the Alzheimer's
# Disease / IGAP / Marioni / Posthuma bullets with new sub-bullets, bolds
# the three study-name bullets, and relocates the _GoBack bookmark into
# the newly-added SusieR sentence.

$d = $word.ActiveDocument

# --- locate the three anchor paragraphs under "Alzheimer's Disease" ---
$igapPara = $null
$marioniPara = $null
$posthumaPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.Trim()
    if ($t -eq "IGAP" -and $igapPara -eq $null) {
        $igapPara = $i
    } elseif ($t -like "Marioni et al*") {
        $marioniPara = $i
    } elseif ($t -like "Posthuma et al*") {
        $posthumaPara = $i
    }
}

# --- move the _GoBack bookmark off of "Alzheimer's Disease" (it will be
#     re-added inside the new SusieR sentence below) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- insert "Study Info" / "PTK2B" / SusieR sentence after IGAP.
#     Done *before* IGAP is bolded so the new paragraphs naturally
#     inherit "not bold" with no explicit <w:b> override. ---
$igap = $d.Paragraphs.Item($igapPara)
$igap.Range.InsertParagraphAfter()
$studyInfo = $d.Paragraphs.Item($igapPara + 1)
$studyInfo.Range.Text = "Study Info"
$studyInfo.Range.ListFormat.ListLevelNumber = 3

$studyInfo.Range.InsertParagraphAfter()
$ptk2bA = $d.Paragraphs.Item($igapPara + 2)
$ptk2bA.Range.Text = "PTK2B"
$ptk2bA.Range.ListFormat.ListLevelNumber = 3

$ptk2bA.Range.InsertParagraphAfter()
$susieA = $d.Paragraphs.Item($igapPara + 3)
$susieA.Range.Text = "SusieR identified rs[#####] as the variant with the top PIP, in contrast to the GWAS where the lead variant was identified as rs[#####]. Specifying that there was more than one causal variant (L>1) prevented the model from converging"
$susieA.Range.ListFormat.ListLevelNumber = 4

$gwasSpot = $susieA.Range.Duplicate
$gwasSpot.Find.Execute("GWAS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gwasSpot.Collapse(1)
$d.Bookmarks.Add("_GoBack", $gwasSpot)

# --- Marioni paragraph number shifted by the 3 new paragraphs above ---
$marioniPara = $marioniPara + 3
$posthumaPara = $posthumaPara + 3

# --- insert meta-analysis note / "PTK2B" / blank line after Marioni,
#     again before Marioni gets bolded. ---
$marioni = $d.Paragraphs.Item($marioniPara)
$marioni.Range.InsertParagraphAfter()
$metaNote = $d.Paragraphs.Item($marioniPara + 1)
$metaNote.Range.Text = "This study was a meta-analysis of a variety of other studies, including IGAP."
$metaNote.Range.ListFormat.ListLevelNumber = 3

$metaNote.Range.InsertParagraphAfter()
$ptk2bB = $d.Paragraphs.Item($marioniPara + 2)
$ptk2bB.Range.Text = "PTK2B"
$ptk2bB.Range.ListFormat.ListLevelNumber = 3

$ptk2bB.Range.InsertParagraphAfter()
$blankB = $d.Paragraphs.Item($marioniPara + 3)
$blankB.Range.ListFormat.ListLevelNumber = 4

# --- Posthuma shifted by the 3 new paragraphs inserted after Marioni ---
$posthumaPara = $posthumaPara + 3

# --- now bold the three study-name bullets ---
$d.Paragraphs.Item($igapPara).Range.Bold = $true
$d.Paragraphs.Item($marioniPara).Range.Bold = $true
$d.Paragraphs.Item($posthumaPara).Range.Bold = $true
